$d = $word.ActiveDocument

# Locate the paragraph that starts "He pulled open the kitchen drawer..." and
# expand the (collapsed) found range out to the whole enclosing paragraph so we
# can swap its contents for the new, expanded passage below.
$rng = $d.Content
$found = $rng.Find.Execute("He pulled open the kitchen drawer and found a notepad", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate target paragraph"
}
$rng.Expand(4) | Out-Null

# Replace that whole paragraph with three paragraphs of fresh OOXML: the
# original opening plus "...every meaningful human connection..." and the new
# material about Shelton's phone call, followed by two new paragraphs
# describing the pilot interrupting him.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">He </w:t></w:r><w:r><w:t>pulled open the kitchen drawer and found a notepad</w:t></w:r><w:r><w:t xml:space="preserve">. Scratched into the cardboard at the back of the pad were </w:t></w:r><w:r><w:t xml:space="preserve">a set of names and numbers in his father’s handwriting. </w:t></w:r><w:r><w:t xml:space="preserve">Less than </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">10 names that represented every </w:t></w:r><w:r><w:t xml:space="preserve">meaningful </w:t></w:r><w:r><w:t>human connection he’d made during his life.</w:t></w:r><w:r><w:t xml:space="preserve"> Shelton’s old number was on there – the Pikeville place – and Abby’s number – first, her number from outside Chicago where she settled after school; that one was crossed out with big springy loops and next to it was written the number Shelton had just dialed. Shelton found the number for Abby’s mother-in-law beside the word “Joyce” at the bottom of the list and dialed. She lived across town. Abby and Mark were always visiting, dropping off the kids. She would know how to reach Abby. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Receiver pressed to his left ear and his palm pressed hard against the right, Shelton heard nothing but the hollow ringing of an unanswered line. But he became aware of the cabin door opening behind him. It was a feeling – the cold blowing through the cabin, the shifting light as shadows were let in. </w:t></w:r><w:r><w:t>A hand tapped his shoulder and Shelton turned.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">“Why don’t you make us some tea?” The pilot slipped the receiver from Shelton’s slack grip and set it back in its cradle. </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml) | Out-Null

Write-Host "Paragraphs now:" $d.Paragraphs.Count
